$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin) - D2 stays text-like (has multiple dots), no apostrophe needed
$ws.Range("D2").Value = "65.310.28"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3 (Ethereum) - D3 stays text-like
$ws.Range("D3").Value = "2.950.13"
$ws.Range("E3").Value = "  -2.48%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 (BNB) - D5 is a plain number, force text to avoid numeric conversion
$ws.Range("D5").Value = "'572.05"
$ws.Range("E5").Value = "  -2.36%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "'161.19"
$ws.Range("E6").Value = "  +0.45%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 (XRP)
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  -0.45%  "

# Row 9 (LidoStakedEther) - D9 stays text-like
$ws.Range("D9").Value = "2.947.95"
$ws.Range("E9").Value = "  -2.42%  "

# Row 10 (Toncoin)
$ws.Range("D10").Value = "'6.74"
$ws.Range("E10").Value = "  -3.37%  "

# Row 11 (Dogecoin)
$ws.Range("E11").Value = "  -2.97%  "

# Row 12 (Cardano)
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13 (ShibaInu)
$ws.Range("E13").Value = "  -2.76%  "

# Row 14 (Avalanche)
$ws.Range("D14").Value = "'34.48"
$ws.Range("E14").Value = "  -0.41%  "

# Row 15 (TRON)
$ws.Range("E15").Value = "  -0.98%  "

# Row 16 (WrappedBTC) - D16 stays text-like
$ws.Range("D16").Value = "65.338.96"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17 (WrappedliquidstakedEther2.0) - D17 stays text-like
$ws.Range("D17").Value = "3.438.30"
$ws.Range("E17").Value = "  -2.44%  "

# Row 18 (Polkadot)
$ws.Range("D18").Value = "'7.06"
$ws.Range("E18").Value = "  +1.39%  "

# Row 19 (WrappedEther) - D19 stays text-like
$ws.Range("D19").Value = "2.947.19"
$ws.Range("E19").Value = "  -2.62%  "

# Row 20 (Chainlink)
$ws.Range("D20").Value = "'15.57"
$ws.Range("E20").Value = "  +12.38%  "

# Row 21 (BitcoinCash)
$ws.Range("D21").Value = "'444.64"
$ws.Range("E21").Value = "  -3.37%  "

# Row 22 (Polygon)
$ws.Range("D22").Value = "'0.693"
$ws.Range("E22").Value = "  +1.19%  "

# Row 23 (Uniswap)
$ws.Range("E23").Value = "  -0.89%  "

# Row 24 (Litecoin)
$ws.Range("D24").Value = "'82.63"
$ws.Range("E24").Value = "  +0.25%  "

# Row 25 (Fetch.AI)
$ws.Range("E25").Value = "  +1.00%  "

# Row 26 (InternetComputer(DFINITY))
$ws.Range("D26").Value = "'12.16"
$ws.Range("E26").Value = "  -2.28%  "

# Row 27 (RenderToken)
$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  -5.44%  "

# Row 28 (Dai)
$ws.Range("E28").Value = "  +0.11%  "

# Row 29: NEARProtocol -> ImmutableX
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'2.45"
$ws.Range("E29").Value = "  +2.70%  "

# Row 30: ImmutableX -> NEARProtocol
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'8.04"
$ws.Range("E30").Value = "  -0.42%  "

# Row 31 (PancakeSwap)
$ws.Range("E31").Value = "  -0.69%  "

# Row 32 (PEPE)
$ws.Range("E32").Value = "  -1.98%  "

# Row 33: Hedera -> EthereumClassic
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'27.25"
$ws.Range("E33").Value = "  +1.00%  "

# Row 34: EthereumClassic -> Hedera
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.112"
$ws.Range("E34").Value = "  +1.42%  "

# Row 35 (FirstDigitalUSD)
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.15%  "

# Row 36 (Mantle)
$ws.Range("D36").Value = "'0.974"
$ws.Range("E36").Value = "  -2.16%  "

# Row 37 (Filecoin)
$ws.Range("E37").Value = "  -0.14%  "

# Row 38 (OKB)
$ws.Range("D38").Value = "'49.64"
$ws.Range("E38").Value = "  -0.40%  "

# Row 39 (Arweave)
$ws.Range("D39").Value = "'45.14"
$ws.Range("E39").Value = "  +1.36%  "

# Row 40: Stacks -> TheGraph
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.302"
$ws.Range("E40").Value = "  +0.00%  "

# Row 41: TheGraph -> Kaspa
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.120"
$ws.Range("E41").Value = "  -0.65%  "

# Row 42: Kaspa -> dogwifhat
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.85"
$ws.Range("E42").Value = "  -6.98%  "

# Row 43: dogwifhat -> Stacks
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.94"
$ws.Range("E43").Value = "  -10.20%  "

# Row 44 (Cosmos)
$ws.Range("D44").Value = "'8.57"
$ws.Range("E44").Value = "  +1.55%  "

# Row 45 (Bittensor)
$ws.Range("D45").Value = "'383.59"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46 (VeChain)
$ws.Range("D46").Value = "'0.0352"
$ws.Range("E46").Value = "  -0.73%  "

# Row 47 (Maker) - D47 stays text-like
$ws.Range("D47").Value = "2.697.47"
$ws.Range("E47").Value = "  -3.97%  "

# Row 48 (Monero)
$ws.Range("D48").Value = "'134.51"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49 (USDe)
$ws.Range("E49").Value = "  +0.03%  "

# Row 50 (ThetaToken)
$ws.Range("D50").Value = "'2.19"
$ws.Range("E50").Value = "  +3.78%  "

# Row 51 (InjectiveProtocol)
$ws.Range("D51").Value = "'23.54"
$ws.Range("E51").Value = "  -0.40%  "
